$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("C4").Value = 25
$ws.Range("B5").Value = 0.95
$ws.Range("C5").Value = 1.45

# Update selection to C5
$ws.Range("C5").Select()
